$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.022855043411255
$ws.Range("B1").Value = 1.898676156997681
$ws.Range("C1").Value = 7.745347499847412
$ws.Range("D1").Value = 2.277709007263184
$ws.Range("E1").Value = 0.4974815845489502
